# Resort the worksheets: move "总计" so it comes before "2021-Q3".
# (The underlying data/formatting of each sheet is unchanged - only the
# tab order/sheet names swap position, per the commit "update data with
# resort sheetname".)

$wb = $excel.ActiveWorkbook

# Move "总计" so it sits immediately before "2021-Q3" (i.e. becomes the
# first tab). Re-resolve sheets by name rather than caching references
# across the Move call.
$wb.Worksheets.Item("总计").Move($wb.Worksheets.Item("2021-Q3"))

# Keep "2021-Q3" as the active/selected tab, matching the original
# selection state of the workbook.
$wb.Worksheets.Item("2021-Q3").Activate()
